$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Onion white: clear "Current Quantity" (B2, keep as text) and set "Requested quantity" (C2) to 0
$ws.Range("B2").Value = "'"
$ws.Range("C2").Value = 0

# Row 3 - Onion Red: clear "Current Quantity" (B3, keep as text)
$ws.Range("B3").Value = "'"

# Row 4 - Potato Russel: set "Current Quantity" (B4) to text "1" and "Requested quantity" (C4) to 3
$ws.Range("B4").Value = "'1"
$ws.Range("C4").Value = 3

# Row 21 - Frozen Peas and carrot cut: set B21 to text "2" and C21 to 1
$ws.Range("B21").Value = "'2"
$ws.Range("C21").Value = 1

# Row 39 - Paneer - Not Appel, not Nanak: set B39 to text "1" and C39 to 2
$ws.Range("B39").Value = "'1"
$ws.Range("C39").Value = 2
